$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append: row, dateSerial, colB, colC, colD
$data = @(
    @(329, 44403, 0, 5, 31.30870381966186),
    @(330, 44404, 0, 5, 31.30870381966186),
    @(331, 44405, 0, 5, 31.30870381966186),
    @(332, 44406, 2, 3, 18.78522229179712),
    @(333, 44407, 3, 6, 37.57044458359425),
    @(334, 44408, 0, 5, 31.30870381966186),
    @(335, 44409, 0, 5, 31.30870381966186),
    @(336, 44410, 1, 6, 37.57044458359425),
    @(337, 44411, 2, 8, 50.09392611145898),
    @(338, 44412, 1, 9, 56.35566687539136),
    @(339, 44413, 4, 11, 68.8791484032561),
    @(340, 44414, 0, 8, 50.09392611145898),
    @(341, 44415, 2, 10, 62.61740763932373),
    @(342, 44416, 1, 11, 68.8791484032561),
    @(343, 44417, 1, 11, 68.8791484032561)
)

foreach ($r in $data) {
    $row = $r[0]
    # Copy the formatting (including date style) from the row above, then set the actual values.
    $ws.Cells.Item($row - 1, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}
